# Update simulation results (Results_S6) for rows 2-29, columns B:H
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.13602999999475
$ws.Range("C2").Value = 381
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0.2392189636994589
$ws.Range("G2").Value = 3696.686773354967
$ws.Range("H2").Value = 0.5717560425281184
$ws.Range("B3").Value = 21.85007670999479
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.2821717707289705
$ws.Range("G3").Value = 3777.306414268884
$ws.Range("H3").Value = 0.5784565590828293
$ws.Range("B4").Value = 22.56042067999476
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.3027178665982671
$ws.Range("G4").Value = 3918.703817962509
$ws.Range("H4").Value = 0.575711299654303
$ws.Range("B5").Value = 23.22582387999479
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0.3158736445565552
$ws.Range("G5").Value = 4077.638135965924
$ws.Range("H5").Value = 0.569590118238704
$ws.Range("B6").Value = 23.95771868999476
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0.3200592976877066
$ws.Range("G6").Value = 4259.609038422821
$ws.Range("H6").Value = 0.562439380560274
$ws.Range("B7").Value = 24.76975238999474
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 1.218841982306258
$ws.Range("F7").Value = 0.3659390839324042
$ws.Range("G7").Value = 4461.549043111811
$ws.Range("H7").Value = 0.5551827885482236
$ws.Range("B8").Value = 25.66272002999472
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 14.662803998
$ws.Range("F8").Value = 0.3661127750801957
$ws.Range("G8").Value = 4599.804730061753
$ws.Range("H8").Value = 0.5579089012687328
$ws.Range("B9").Value = 26.74449718999474
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 32.20124931477594
$ws.Range("F9").Value = 0.386180324527109
$ws.Range("G9").Value = 4677.348646816087
$ws.Range("H9").Value = 0.5717875490893748
$ws.Range("B10").Value = 27.10424367999474
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = 82.45031263246121
$ws.Range("F10").Value = 0.3874726101052314
$ws.Range("G10").Value = 4755.59162283329
$ws.Range("H10").Value = 0.5699447267477217
$ws.Range("B11").Value = 27.46370780999474
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 18
$ws.Range("E11").Value = 115.7294124388981
$ws.Range("F11").Value = 0.3809946954220956
$ws.Range("G11").Value = 4877.09837581188
$ws.Range("H11").Value = 0.5631157236073369
$ws.Range("B12").Value = 27.83064783999474
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 31
$ws.Range("E12").Value = 275.9049874085127
$ws.Range("F12").Value = 0.334889149948146
$ws.Range("G12").Value = 5063.836880365056
$ws.Range("H12").Value = 0.5495960572487559
$ws.Range("B13").Value = 28.18640283999475
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 35
$ws.Range("E13").Value = 272.7165734417719
$ws.Range("F13").Value = 0.3236324553593016
$ws.Range("G13").Value = 5237.173940461217
$ws.Range("H13").Value = 0.5381987148113031
$ws.Range("B14").Value = 28.47498796999474
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 32
$ws.Range("E14").Value = 255.4813828723899
$ws.Range("F14").Value = 0.3237395415563139
$ws.Range("G14").Value = 5340.320279327852
$ws.Range("H14").Value = 0.5332074946931589
$ws.Range("B15").Value = 28.52692242179329
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 25
$ws.Range("E15").Value = 277.8570413816718
$ws.Range("F15").Value = 0.3106979675859447
$ws.Range("G15").Value = 5416.090169571141
$ws.Range("H15").Value = 0.5267069330208756
$ws.Range("B16").Value = 28.56714351999474
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 27
$ws.Range("E16").Value = 210.5338510351209
$ws.Range("F16").Value = 0.2854469826585077
$ws.Range("G16").Value = 5507.814321845025
$ws.Range("H16").Value = 0.5186656966029535
$ws.Range("B17").Value = 28.59771566999473
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = 179.6003763493269
$ws.Range("F17").Value = 0.2752346970721594
$ws.Range("G17").Value = 5565.924514290129
$ws.Range("H17").Value = 0.5137999194306724
$ws.Range("B18").Value = 28.59428545999473
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 22
$ws.Range("E18").Value = 169.901567316808
$ws.Range("F18").Value = 0.2753960011405997
$ws.Range("G18").Value = 5688.476045299411
$ws.Range("H18").Value = 0.5026704029741533
$ws.Range("B19").Value = 28.57422669999474
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = 139.7670312623993
$ws.Range("F19").Value = 0.277828563280043
$ws.Range("G19").Value = 5750.192130385776
$ws.Range("H19").Value = 0.4969264687522313
$ws.Range("B20").Value = 28.3714426800693
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 92.03246384220219
$ws.Range("F20").Value = 0.2769517674928316
$ws.Range("G20").Value = 5740.545051165566
$ws.Range("H20").Value = 0.4942290745424729
$ws.Range("B21").Value = 28.15722030999473
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 46.58455201596804
$ws.Range("F21").Value = 0.2579171075090837
$ws.Range("G21").Value = 5727.843918854681
$ws.Range("H21").Value = 0.4915849787266017
$ws.Range("B22").Value = 27.92259855054463
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 28.5671627569966
$ws.Range("F22").Value = 0.2481569928825704
$ws.Range("G22").Value = 5730.626338624557
$ws.Range("H22").Value = 0.4872521239492733
$ws.Range("B23").Value = 27.65148095054423
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 8.714427069735013
$ws.Range("F23").Value = 0.2475440543536748
$ws.Range("G23").Value = 5673.901022849678
$ws.Range("H23").Value = 0.4873451411857104
$ws.Range("B24").Value = 27.35481787054649
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 2.457059041737089
$ws.Range("F24").Value = 0.247545411270181
$ws.Range("G24").Value = 5620.195477581339
$ws.Range("H24").Value = 0.4867236020466442
$ws.Range("B25").Value = 26.97910179054833
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 53.35948924898733
$ws.Range("F25").Value = 0.2439679545767239
$ws.Range("G25").Value = 5520.897437838777
$ws.Range("H25").Value = 0.4886723960065019
$ws.Range("B26").Value = 26.59093877067457
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 123.899344982822
$ws.Range("F26").Value = 0.2406115871854294
$ws.Range("G26").Value = 5416.471797912552
$ws.Range("H26").Value = 0.4909273003308615
$ws.Range("B27").Value = 26.22162253097315
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 20
$ws.Range("E27").Value = 218.6215476182707
$ws.Range("F27").Value = 0.2352180149906909
$ws.Range("G27").Value = 5396.19152931093
$ws.Range("H27").Value = 0.4859283142294531
$ws.Range("B28").Value = 25.86659505066243
$ws.Range("C28").Value = 14
$ws.Range("D28").Value = 32
$ws.Range("E28").Value = 298.2008961105664
$ws.Range("F28").Value = 0.2246198124795308
$ws.Range("G28").Value = 5353.65965256259
$ws.Range("H28").Value = 0.4831572555846185
$ws.Range("B29").Value = 25.5086743001779
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 38
$ws.Range("E29").Value = 270.5505140766452
$ws.Range("F29").Value = 0.2252550489982136
$ws.Range("G29").Value = 5375.70732750707
$ws.Range("H29").Value = 0.474517542457195
